# Weekly refresh of Fruta/Hortaliza price data for Arandano (blue):
# a new week (row 37) is inserted at the top of the date-ordered history
# and the existing rows 37-50 each shift down by one, taking on the
# prior rows values for Fecha (D), Calidad (L), Volumen (M),
# Precio minimo/maximo/promedio (N/O/P), Origen (R) and Precio $/Kg (S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37
$ws.Range("D37").Value = 44624
$ws.Range("M37").Value = 380
$ws.Range("N37").Value = 4000
$ws.Range("O37").Value = 4000
$ws.Range("P37").Value = 4000
$ws.Range("R37").Value = 'Provincia de Linares'
$ws.Range("S37").Value = 2000

# Row 38
$ws.Range("D38").Value = 44607
$ws.Range("M38").Value = 400
$ws.Range("N38").Value = 3400
$ws.Range("O38").Value = 3600
$ws.Range("P38").Value = 3510
$ws.Range("R38").Value = 'Provincia de Curicó'
$ws.Range("S38").Value = 1755

# Row 39
$ws.Range("D39").Value = 44580
$ws.Range("M39").Value = 410
$ws.Range("N39").Value = 4000
$ws.Range("O39").Value = 4000
$ws.Range("P39").Value = 4000
$ws.Range("S39").Value = 2000

# Row 40
$ws.Range("D40").Value = 44189
$ws.Range("M40").Value = 180
$ws.Range("R40").Value = 'Región de O''Higgins'

# Row 41
$ws.Range("D41").Value = 44187
$ws.Range("M41").Value = 250
$ws.Range("N41").Value = 3000
$ws.Range("O41").Value = 3000
$ws.Range("P41").Value = 3000
$ws.Range("R41").Value = 'Provincia de Linares'
$ws.Range("S41").Value = 1500

# Row 42
$ws.Range("D42").Value = 44449
$ws.Range("M42").Value = 65
$ws.Range("N42").Value = 16000
$ws.Range("O42").Value = 16000
$ws.Range("P42").Value = 16000
$ws.Range("R42").Value = 'Provincia del Elquí'
$ws.Range("S42").Value = 8000

# Row 43
$ws.Range("L43").Value = 'Primera'
$ws.Range("M43").Value = 100
$ws.Range("N43").Value = 3600
$ws.Range("O43").Value = 3600
$ws.Range("P43").Value = 3600
$ws.Range("S43").Value = 1800

# Row 44
$ws.Range("D44").Value = 44235
$ws.Range("L44").Value = 'Segunda'
$ws.Range("M44").Value = 150
$ws.Range("N44").Value = 2600
$ws.Range("O44").Value = 2600
$ws.Range("P44").Value = 2600
$ws.Range("R44").Value = 'Provincia de Curicó'
$ws.Range("S44").Value = 1300

# Row 45
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 450
$ws.Range("N45").Value = 2800
$ws.Range("O45").Value = 3000
$ws.Range("P45").Value = 2911
$ws.Range("S45").Value = 1456

# Row 46
$ws.Range("D46").Value = 44209
$ws.Range("L46").Value = 'Segunda'
$ws.Range("M46").Value = 500
$ws.Range("N46").Value = 2500
$ws.Range("O46").Value = 2600
$ws.Range("P46").Value = 2556
$ws.Range("S46").Value = 1278

# Row 47
$ws.Range("L47").Value = 'Primera'
$ws.Range("M47").Value = 200
$ws.Range("N47").Value = 2800
$ws.Range("O47").Value = 3000
$ws.Range("P47").Value = 2920
$ws.Range("S47").Value = 1460

# Row 48
$ws.Range("D48").Value = 44250
$ws.Range("L48").Value = 'Segunda'
$ws.Range("M48").Value = 300
$ws.Range("N48").Value = 2400
$ws.Range("O48").Value = 2500
$ws.Range("P48").Value = 2453
$ws.Range("R48").Value = 'Provincia de Linares'
$ws.Range("S48").Value = 1226

# Row 49
$ws.Range("D49").Value = 44181
$ws.Range("M49").Value = 220
$ws.Range("N49").Value = 3600
$ws.Range("O49").Value = 4000
$ws.Range("P49").Value = 3782
$ws.Range("R49").Value = 'Región de O''Higgins'
$ws.Range("S49").Value = 1891

# Row 50
$ws.Range("D50").Value = 44614
$ws.Range("M50").Value = 400
$ws.Range("N50").Value = 3800
$ws.Range("O50").Value = 3800
$ws.Range("P50").Value = 3800
$ws.Range("R50").Value = 'Provincia de Curicó'
$ws.Range("S50").Value = 1900
